$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.99999999825042629
$ws.Cells.Item(2, 1).Value = 0.99767544249126128
$ws.Cells.Item(3, 1).Value = 0.98884576412874292
$ws.Cells.Item(4, 1).Value = 0.99167482224131409
$ws.Cells.Item(5, 1).Value = 0.98135394293647726
$ws.Cells.Item(6, 1).Value = 0.95656583680300833
$ws.Cells.Item(7, 1).Value = 0.95067909061784195
$ws.Cells.Item(8, 1).Value = 0.94280897649427275
$ws.Cells.Item(9, 1).Value = 0.93313230188897278
$ws.Cells.Item(10, 1).Value = 0.92437096503852345
$ws.Cells.Item(11, 1).Value = 0.92311409150469426
$ws.Cells.Item(12, 1).Value = 0.9209791105009647
$ws.Cells.Item(13, 1).Value = 0.9126284905564821
$ws.Cells.Item(14, 1).Value = 0.90846089529060214
$ws.Cells.Item(15, 1).Value = 0.9058692610243182
$ws.Cells.Item(16, 1).Value = 0.90336267525205649
$ws.Cells.Item(17, 1).Value = 0.89965461864970075
$ws.Cells.Item(18, 1).Value = 0.89854567575415856
$ws.Cells.Item(19, 1).Value = 0.99662217206931869
$ws.Cells.Item(20, 1).Value = 0.98950495039613018
$ws.Cells.Item(21, 1).Value = 0.98810642554132144
$ws.Cells.Item(22, 1).Value = 0.98684191425524237
$ws.Cells.Item(23, 1).Value = 0.97714300071137039
$ws.Cells.Item(24, 1).Value = 0.95555101942132925
$ws.Cells.Item(25, 1).Value = 0.94909389241395792
$ws.Cells.Item(26, 1).Value = 0.9545807252101578
$ws.Cells.Item(27, 1).Value = 0.95206318416673796
$ws.Cells.Item(28, 1).Value = 0.94153148982732893
$ws.Cells.Item(29, 1).Value = 0.93439278729436159
$ws.Cells.Item(30, 1).Value = 0.93183035131619918
$ws.Cells.Item(31, 1).Value = 0.93456775331157682
$ws.Cells.Item(32, 1).Value = 0.93651733148702543
$ws.Cells.Item(33, 1).Value = 0.93599730124377523
